# Esquema definitivo del TFG
#
# The heading "4. Plan de pruebas e implementación (2 - 4 // 5+)" and the
# sub-heading "4.1. Plan de implementación [PROYECTO en la empresa] (1 - 2)"
# both have the word "implementación" corrected to "implantación" (the
# middle "eme" is replaced by a single "a"). Word records the caret's last
# edit position with the hidden "_GoBack" bookmark, which ends up right
# after the newly typed "a" in the first heading (and is removed from its
# previous location at the very end of the document).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Heading 4: "4. Plan de pruebas e implementación (2 - 4 // 5+)"
#                                          ^^^ -> "a"
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("pruebas e implementación (2", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found1) { throw "Could not find heading 4 target text" }

$emeStart1 = $rng1.Start + 14   # offset of "eme" inside "pruebas e implementación (2"
$emeEnd1   = $emeStart1 + 3

$emeRng1 = $d.Range($emeStart1, $emeEnd1)
if ($emeRng1.Text -ne "eme") { throw "Unexpected text at heading 4 target: $($emeRng1.Text)" }
$emeRng1.Text = "a"

# Force a run boundary right before the new "a" so it becomes its own run,
# matching the fine-grained run split Word produced.
$split1 = $d.Range($emeStart1, $emeStart1)
$d.Bookmarks.Add("zzSplitHeading4", $split1) | Out-Null

# The caret (and therefore "_GoBack") ends up right after the new "a".
$bmPos1 = $emeStart1 + 1
$bmRange1 = $d.Range($bmPos1, $bmPos1)
$d.Bookmarks.Add("_GoBack", $bmRange1) | Out-Null

# Drop the scaffolding bookmark now that the split is baked into the runs.
$d.Bookmarks("zzSplitHeading4").Delete()

# ---------------------------------------------------------------------
# Heading 4.1: "4.1. Plan de implementación [PROYECTO en la empresa] (1 - 2)"
#                                 ^^^ -> "a"
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("implementación [PROYECTO", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find heading 4.1 target text" }

$emeStart2 = $rng2.Start + 4    # offset of "eme" inside "implementación [PROYECTO"
$emeEnd2   = $emeStart2 + 3

# Split the run before "eme" first (while it still reads "eme") so the tab
# that precedes this text in the same run is preserved on the left piece.
$splitBefore2 = $d.Range($emeStart2, $emeStart2)
$d.Bookmarks.Add("zzSplitHeading41a", $splitBefore2) | Out-Null
$d.Bookmarks("zzSplitHeading41a").Delete()

$emeRng2 = $d.Range($emeStart2, $emeEnd2)
if ($emeRng2.Text -ne "eme") { throw "Unexpected text at heading 4.1 target: $($emeRng2.Text)" }
$emeRng2.Text = "a"

# Split again right after the new "a" so it lands in its own run too.
$splitAfter2 = $d.Range($emeStart2 + 1, $emeStart2 + 1)
$d.Bookmarks.Add("zzSplitHeading41b", $splitAfter2) | Out-Null
$d.Bookmarks("zzSplitHeading41b").Delete()
